$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update values in column E ("F" header) for the affected rows
$ws.Range("E3").Value = 12.968
$ws.Range("E21").Value = 13.492
$ws.Range("E23").Value = 13.121
$ws.Range("E25").Value = 12.659
